$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I20").Value = -0.475090227869142
$ws.Range("J20").Value = 0.3390417815583261
$ws.Range("K20").Value = 0.2941700511571806
$ws.Range("L20").Value = 2.206662329477553
